$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (values like "0.513" or "2.18" would otherwise be auto-converted to numbers)
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "44.563.66"
$ws.Range("E2").Value = "  +3.89%  "

$ws.Range("D3").Value = "2.430.74"
$ws.Range("E3").Value = "  +2.91%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "313.10"
$ws.Range("E5").Value = "  +3.78%  "

$ws.Range("D6").Value = "101.82"
$ws.Range("E6").Value = "  +6.87%  "

$ws.Range("D7").Value = "0.513"
$ws.Range("E7").Value = "  +1.82%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "0.513"
$ws.Range("E9").Value = "  +5.70%  "

$ws.Range("D10").Value = "35.29"
$ws.Range("E10").Value = "  +4.26%  "

$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("D12").Value = "0.125"
$ws.Range("E12").Value = "  +1.32%  "

$ws.Range("D13").Value = "18.77"
$ws.Range("E13").Value = "  +2.89%  "

$ws.Range("E14").Value = "  +3.56%  "

$ws.Range("D15").Value = "2.810.01"
$ws.Range("E15").Value = "  +2.78%  "

$ws.Range("D16").Value = "2.418.54"
$ws.Range("E16").Value = "  +2.78%  "

$ws.Range("E17").Value = "  +5.46%  "

$ws.Range("D18").Value = "44.497.22"
$ws.Range("E18").Value = "  +3.87%  "

$ws.Range("D19").Value = "12.43"
$ws.Range("E19").Value = "  +3.34%  "

$ws.Range("E20").Value = "  +2.26%  "

$ws.Range("E21").Value = "  +2.55%  "

$ws.Range("D22").Value = "68.96"
$ws.Range("E22").Value = "  +1.66%  "

$ws.Range("D23").Value = "241.23"
$ws.Range("E23").Value = "  +2.66%  "

$ws.Range("E24").Value = "  +4.36%  "

$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  +2.54%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").Value = "25.20"
$ws.Range("E27").Value = "  +2.32%  "

$ws.Range("E28").Value = "  -4.17%  "

$ws.Range("D29").Value = "9.65"
$ws.Range("E29").Value = "  +4.44%  "

$ws.Range("D30").Value = "33.29"
$ws.Range("E30").Value = "  +5.88%  "

$ws.Range("D31").Value = "48.47"
$ws.Range("E31").Value = "  +1.36%  "

$ws.Range("D32").Value = "0.123"
$ws.Range("E32").Value = "  +17.02%  "

$ws.Range("D33").Value = "19.51"
$ws.Range("E33").Value = "  +12.92%  "

$ws.Range("E34").Value = "  +3.54%  "

$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("E36").Value = "  +5.87%  "

$ws.Range("D37").Value = "1.90"
$ws.Range("E37").Value = "  +2.92%  "

$ws.Range("D38").Value = "4.54"
$ws.Range("E38").Value = "  +4.26%  "

$ws.Range("E39").Value = "  +4.48%  "

$ws.Range("D40").Value = "127.24"
$ws.Range("E40").Value = "  +5.37%  "

$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "22.06"
$ws.Range("E41").Value = "  +2.99%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.109"
$ws.Range("E42").Value = "  +1.05%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "2.18"
$ws.Range("E43").Value = "  -5.38%  "

$ws.Range("D44").Value = "0.0289"
$ws.Range("E44").Value = "  +3.68%  "

$ws.Range("D45").Value = "1.947.10"
$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("E46").Value = "  +2.31%  "

$ws.Range("E47").Value = "  +8.72%  "

$ws.Range("E48").Value = "  +6.81%  "

$ws.Range("D49").Value = "1.68"
$ws.Range("E49").Value = "  +11.58%  "

$ws.Range("D50").Value = "53.54"
$ws.Range("E50").Value = "  +3.80%  "

$ws.Range("D51").Value = "73.88"
$ws.Range("E51").Value = "  +2.71%  "

# Reset style index so no stray cell-level style is introduced by the
# NumberFormat change above (keeps styles.xml untouched for these cells)
$priceVolRange.Style = "Normal"
